# Generate Report for Handback
#
# The localization-status report is refreshed: every row that was
# "Ready for handoff" has now been handed back and is in sync with en-US,
# so the Status column is updated everywhere it appears (Overview summary
# sheet plus each per-locale detail sheet), and each per-locale sheet grows
# two new columns recording the file that was handed back (same file as the
# original source/handoff columns) and the datetime of that handback.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: refresh the Status shown for both locale columns on
# every row that previously said "Ready for handoff".
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
for ($r = 2; $r -le 3; $r++) {
    foreach ($col in @("B", "C")) {
        $cell = $overview.Range($col + $r)
        if ($cell.Value2 -eq $statusOld) {
            $cell.Value = $statusNew
        }
    }
}

# ---------------------------------------------------------------------
# Per-locale detail sheets: zh-cn and de-de share the same shape, so loop
# over them with the locale-specific handback timestamp.
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-03-11 01:31:32" },
    @{ Name = "de-de"; HandbackTime = "2016-03-11 01:31:53" }
)

# Hyperlink-colored-underline style, matching the look already used for
# the Source File Name / Latest Handoff File link columns (A and C).
$linkColor = 15570276  # OLE BGR for RGB 64,95,ED (#6495ED)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    for ($r = 2; $r -le 3; $r++) {
        $statusCell = $ws.Range("B" + $r)
        if ($statusCell.Value2 -eq $statusOld) {
            $statusCell.Value = $statusNew
        }

        $sourceFile = $ws.Range("A" + $r).Value2
        $sourceDisplay = $ws.Range("A" + $r).Text
        $handoffFile = $ws.Range("C" + $r).Value2
        $handoffDisplay = $ws.Range("C" + $r).Text

        # Latest Target File (E) mirrors the Source File Name (A) — the
        # handed-back file targets the same source doc.
        $eCell = $ws.Range("E" + $r)
        $eCell.Value = $sourceFile
        $eCell.Font.Underline = $true
        $eCell.Font.Color = $linkColor

        # Latest Handback File (F) mirrors the Latest Handoff File (C) —
        # same translated artifact, now handed back.
        $fCell = $ws.Range("F" + $r)
        $fCell.Value = $handoffFile
        $fCell.Font.Underline = $true
        $fCell.Font.Color = $linkColor

        # Latest Handback DateTime (G) — stamp of this handback.
        $ws.Range("G" + $r).Value = $locale.HandbackTime

        # Give the two new file-name columns real hyperlinks, the same
        # way the existing Source File Name / Latest Handoff File columns
        # are linked.
        $ws.Hyperlinks.Add($eCell, "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $locale.Name + "/" + $sourceDisplay, "", "", $sourceDisplay) | Out-Null
        $ws.Hyperlinks.Add($fCell, "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $locale.Name + "/" + $handoffDisplay, "", "", $handoffDisplay) | Out-Null
    }
}
